$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.059823021316717
$ws.Range("D2").Value = 1.065182748514129
$ws.Range("E2").Value = 1.055647459051967
$ws.Range("F2").Value = 1.074988047093087
$ws.Range("I2").Value = 1.051889597260631
$ws.Range("J2").Value = 1.064807055956093
$ws.Range("K2").Value = 1.06789667691699
$ws.Range("L2").Value = 1.05838734145741
$ws.Range("M2").Value = 1.077675805665224
$ws.Range("N2").Value = 1.066319204110627
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.061049478872867
$ws.Range("D3").Value = 1.066182766418782
$ws.Range("E3").Value = 1.056703305008444
$ws.Range("F3").Value = 1.076136026717236
$ws.Range("I3").Value = 1.052287913127249
$ws.Range("J3").Value = 1.0656855974847
$ws.Range("K3").Value = 1.068711582409388
$ws.Range("L3").Value = 1.059256052219
$ws.Range("M3").Value = 1.078640203635123
$ws.Range("N3").Value = 1.067198993268975
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.061842857764436
$ws.Range("D4").Value = 1.066829627247504
$ws.Range("E4").Value = 1.057386536400972
$ws.Range("F4").Value = 1.076878971464871
$ws.Range("I4").Value = 1.052544338318232
$ws.Range("J4").Value = 1.066253310833571
$ws.Range("K4").Value = 1.06923804820781
$ws.Range("L4").Value = 1.059817593257679
$ws.Range("M4").Value = 1.079263771061842
$ws.Range("N4").Value = 1.067767512835919
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.062176343239066
$ws.Range("D5").Value = 1.067101516327273
$ws.Range("E5").Value = 1.057673774664834
$ws.Range("F5").Value = 1.077191336488313
$ws.Range("I5").Value = 1.052651825940262
$ws.Range("J5").Value = 1.066491796249156
$ws.Range("K5").Value = 1.069459176034576
$ws.Range("L5").Value = 1.060053528903621
$ws.Range("M5").Value = 1.079525809072336
$ws.Range("N5").Value = 1.068006336928149
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.062232333982593
$ws.Range("D6").Value = 1.067147164722782
$ws.Range("E6").Value = 1.057722003718152
$ws.Range("F6").Value = 1.07724378580268
$ws.Range("I6").Value = 1.052669855214697
$ws.Range("J6").Value = 1.066531828401442
$ws.Range("K6").Value = 1.069496292762423
$ws.Range("L6").Value = 1.060093135577121
$ws.Range("M6").Value = 1.079569800001412
$ws.Range("N6").Value = 1.068046425930682
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.061847314011683
$ws.Range("D7").Value = 1.066833260442815
$ws.Range("E7").Value = 1.057390374460897
$ws.Range("F7").Value = 1.076883145178228
$ws.Range("I7").Value = 1.052545775804396
$ws.Range("J7").Value = 1.066256498198292
$ws.Range("K7").Value = 1.069241003707062
$ws.Range("L7").Value = 1.05982074637534
$ws.Range("M7").Value = 1.079267272856499
$ws.Range("N7").Value = 1.067770704727063
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.06023755534645
$ws.Range("D8").Value = 1.065520754842465
$ws.Range("E8").Value = 1.056004281831189
$ws.Range("F8").Value = 1.075375986919335
$ws.Range("I8").Value = 1.052024481905407
$ws.Range("J8").Value = 1.065104121565526
$ws.Range("K8").Value = 1.06817225098094
$ws.Range("L8").Value = 1.05868104552419
$ws.Range("M8").Value = 1.078001824849673
$ws.Range("N8").Value = 1.066616691587299
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.057399174308826
$ws.Range("D9").Value = 1.063206245350035
$ws.Range("E9").Value = 1.05356198342081
$ws.Range("F9").Value = 1.072721091083212
$ws.Range("I9").Value = 1.05109582501577
$ws.Range("J9").Value = 1.063067603256284
$ws.Range("K9").Value = 1.066282554346301
$ws.Range("L9").Value = 1.056668316292494
$ws.Range("M9").Value = 1.075768356105303
$ws.Range("N9").Value = 1.064577281188466
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.055505599555742
$ws.Range("D10").Value = 1.061662032573637
$ws.Range("E10").Value = 1.051933837223691
$ws.Range("F10").Value = 1.070951700778387
$ws.Range("I10").Value = 1.050469915390834
$ws.Range("J10").Value = 1.06170589767089
$ws.Range("K10").Value = 1.065018381467572
$ws.Range("L10").Value = 1.055323455644023
$ws.Range("M10").Value = 1.074276899276026
$ws.Range("N10").Value = 1.063213641824951
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.054685319328384
$ws.Range("D11").Value = 1.060993071330277
$ws.Range("E11").Value = 1.051228829359015
$ws.Range("F11").Value = 1.070185646665942
$ws.Range("I11").Value = 1.0501972674996
$ws.Range("J11").Value = 1.0611152925378
$ws.Range("K11").Value = 1.064469928848947
$ws.Range("L11").Value = 1.054740380264923
$ws.Range("M11").Value = 1.073630478983468
$ws.Range("N11").Value = 1.062622197964829
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.054380575641063
$ws.Range("D12").Value = 1.060744541928546
$ws.Range("E12").Value = 1.050966955174142
$ws.Range("F12").Value = 1.069901114077054
$ws.Range("I12").Value = 1.05009574905055
$ws.Range("J12").Value = 1.060895767101398
$ws.Range("K12").Value = 1.064266048954084
$ws.Range("L12").Value = 1.054523687308436
$ws.Range("M12").Value = 1.073390276691558
$ws.Range("N12").Value = 1.062402360777127
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.054445946729897
$ws.Range("D13").Value = 1.060797854495199
$ws.Range("E13").Value = 1.051023128228284
$ws.Range("F13").Value = 1.069962146673753
$ws.Range("I13").Value = 1.050117536196161
$ws.Range("J13").Value = 1.060942862792183
$ws.Range("K13").Value = 1.064309789152781
$ws.Range("L13").Value = 1.054570173809232
$ws.Range("M13").Value = 1.073441805130471
$ws.Range("N13").Value = 1.062449523349195
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.054660130258515
$ws.Range("D14").Value = 1.060972528787087
$ws.Range("E14").Value = 1.051207182823921
$ws.Range("F14").Value = 1.070162126836375
$ws.Range("I14").Value = 1.050188880952864
$ws.Range("J14").Value = 1.0610971495327
$ws.Range("K14").Value = 1.064453079348276
$ws.Range("L14").Value = 1.054722470666406
$ws.Range("M14").Value = 1.073610625692764
$ws.Range("N14").Value = 1.06260402919458
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.05479208842418
$ws.Range("D15").Value = 1.061080145062026
$ws.Range("E15").Value = 1.051320584512196
$ws.Range("F15").Value = 1.070285343025931
$ws.Range("I15").Value = 1.050232806340058
$ws.Range("J15").Value = 1.061192190976699
$ws.Range("K15").Value = 1.064541343903342
$ws.Range("L15").Value = 1.054816290799239
$ws.Range("M15").Value = 1.073714629245226
$ws.Range("N15").Value = 1.062699205608332
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.055560031266254
$ws.Range("D16").Value = 1.061706422750935
$ws.Range("E16").Value = 1.051980625947671
$ws.Range("F16").Value = 1.07100254331733
$ws.Range("I16").Value = 1.050487975833877
$ws.Range("J16").Value = 1.061745073461655
$ws.Range("K16").Value = 1.065054758071468
$ws.Range("L16").Value = 1.055362136692049
$ws.Range("M16").Value = 1.07431978709618
$ws.Range("N16").Value = 1.063252873249833
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.056041646071894
$ws.Range("D17").Value = 1.062099187111381
$ws.Range("E17").Value = 1.052394648887167
$ws.Range("F17").Value = 1.071452450751064
$ws.Range("I17").Value = 1.050647601348626
$ws.Range("J17").Value = 1.062091619229029
$ws.Range("K17").Value = 1.065376525191279
$ws.Range("L17").Value = 1.055704331837512
$ws.Range("M17").Value = 1.074699222406238
$ws.Range("N17").Value = 1.063599911151946
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.056322530254894
$ws.Range("D18").Value = 1.062328250393509
$ws.Range("E18").Value = 1.052636140611726
$ws.Range("F18").Value = 1.071714884243594
$ws.Range("I18").Value = 1.050740551484098
$ws.Range("J18").Value = 1.062293659288613
$ws.Range("K18").Value = 1.065564104572151
$ws.Range("L18").Value = 1.05590385699132
$ws.Range("M18").Value = 1.074920481878621
$ws.Range("N18").Value = 1.063802238131588
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.056418298881871
$ws.Range("D19").Value = 1.062406350084131
$ws.Range("E19").Value = 1.052718482946467
$ws.Range("F19").Value = 1.071804369041601
$ws.Range("I19").Value = 1.050772218517461
$ws.Range("J19").Value = 1.062362533791213
$ws.Range("K19").Value = 1.065628047042488
$ws.Range("L19").Value = 1.055971877827722
$ws.Range("M19").Value = 1.074995915722666
$ws.Range("N19").Value = 1.063871210443882
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.055989976854516
$ws.Range("D20").Value = 1.062057050276751
$ws.Range("E20").Value = 1.05235022824317
$ws.Range("F20").Value = 1.071404178908992
$ws.Range("I20").Value = 1.050630491258547
$ws.Range("J20").Value = 1.062054447905702
$ws.Range("K20").Value = 1.065342013178534
$ws.Range("L20").Value = 1.055667624941538
$ws.Range("M20").Value = 1.074658518673769
$ws.Range("N20").Value = 1.063562687041076
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.054597060088715
$ws.Range("D21").Value = 1.060921092880985
$ws.Range("E21").Value = 1.051152983423523
$ws.Range("F21").Value = 1.070103237288329
$ws.Range("I21").Value = 1.050167878466308
$ws.Range("J21").Value = 1.061051720041684
$ws.Range("K21").Value = 1.064410888389174
$ws.Range("L21").Value = 1.054677626162505
$ws.Range("M21").Value = 1.073560914823969
$ws.Range("N21").Value = 1.062558535188476
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.053720958058979
$ws.Range("D22").Value = 1.060206596431725
$ws.Range("E22").Value = 1.050400210766407
$ws.Range("F22").Value = 1.069285363988763
$ws.Range("I22").Value = 1.049875598250481
$ws.Range("J22").Value = 1.060420405924963
$ws.Range("K22").Value = 1.063824525920114
$ws.Range("L22").Value = 1.054054521346779
$ws.Range("M22").Value = 1.072870268866831
$ws.Range("N22").Value = 1.061926324533297
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.054185427495554
$ws.Range("D23").Value = 1.060585391015223
$ws.Range("E23").Value = 1.050799271904971
$ws.Range("F23").Value = 1.069718927001379
$ws.Range("I23").Value = 1.050030676110769
$ws.Range("J23").Value = 1.06075515944292
$ws.Range("K23").Value = 1.064135456099531
$ws.Range("L23").Value = 1.054384903373911
$ws.Range("M23").Value = 1.073236444908902
$ws.Range("N23").Value = 1.062261553439646
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.056013324053568
$ws.Range("D24").Value = 1.062076090189897
$ws.Range("E24").Value = 1.052370300022866
$ws.Range("F24").Value = 1.0714259908455
$ws.Range("I24").Value = 1.050638223056059
$ws.Range("J24").Value = 1.06207124431789
$ws.Range("K24").Value = 1.06535760798658
$ws.Range("L24").Value = 1.055684211428752
$ws.Range("M24").Value = 1.07467691111994
$ws.Range("N24").Value = 1.063579507306096
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.058133187565807
$ws.Range("D25").Value = 1.063804809852723
$ws.Range("E25").Value = 1.054193362087543
$ws.Range("F25").Value = 1.073407344832959
$ws.Range("I25").Value = 1.051337102003237
$ws.Range("J25").Value = 1.063594796575069
$ws.Range("K25").Value = 1.066771852947434
$ws.Range("L25").Value = 1.057189186499579
$ws.Range("M25").Value = 1.076346193204516
$ws.Range("N25").Value = 1.065105223182234
